$d = $word.ActiveDocument

# Locate the end of the "content providers" paragraph (just before the
# existing _GoBack bookmark) and anchor a Range there.
$rng = $d.Content
$rng.Find.Execute("content providers", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# Remove the old _GoBack bookmark - it will be re-created at the end of the
# newly added content, matching where Word leaves it after typing.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Six blank paragraphs.
for ($i = 0; $i -lt 6; $i++) {
    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $rng.MoveStart(1, 1) | Out-Null
}

# "Android Build files - 9.3"
$rng.InsertAfter("Android Build files " + [char]0x2013 + " 9.3")
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null

# DEX/ART/ELF paragraph
$rng.InsertAfter("When build the project is converted into DEX files. When loaded onto the device ART (Android Runtime) does ahead of time compilation which translates the bytecode to native architecture instructions in the ELF (Executable and Linkable Format) format. When the app is launch after that initial load the ELF version is run, increasing app performance. ")
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null

# JIT paragraph
$rng.InsertAfter("Note that in earlier versions of Android JIT (just in time) compilation was used where the bytecode was translated by the VM every time the app was run.")
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null

# One blank paragraph.
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null

# Final paragraph: two runs.
$rng.InsertAfter("Summarize all of chapter 10")
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null
$rng.InsertAfter(" " + [char]0x2013 + " anatomy of and android app (all core parts)")
$rng.Collapse(0)

# Re-create the _GoBack bookmark at the new end-of-document insertion point.
$d.Bookmarks.Add("_GoBack", $rng)
